$d = $word.ActiveDocument
Write-Output ("ParaCount=" + $d.Paragraphs.Count)
Write-Output ("Para1=" + $d.Paragraphs.Item(1).Range.Text)
Write-Output ("ParaLast=" + $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text)
